# Add a second batch of simulation results next to the existing one:
# columns V:AO (20 more columns) get the same angle sequence in row 2
# and fresh 0/1 "collision" results in row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (angles 0.1*pi .. 2*pi, 20 steps) - identical sequence to B2:U2
$row2Values = @(
    0.31415926535897898,
    0.62831853071795896,
    0.94247779607693805,
    1.2566370614359199,
    1.5707963267949001,
    1.8849555921538801,
    2.1991148575128601,
    2.5132741228718301,
    2.8274333882308098,
    3.14159265358979,
    3.4557519189487702,
    3.76991118430775,
    4.0840704496667302,
    4.3982297150257104,
    4.7123889803846897,
    5.0265482457436699,
    5.3407075111026501,
    5.6548667764616303,
    5.9690260418206096,
    6.2831853071795898
)

# Row 3: newly collected 0/1 results for the new batch
$row3Values = @(1, 1, 1, 1, 1, 1, 0, 1, 1, 0, 1, 1, 0, 1, 1, 1, 0, 1, 1, 1)

$n = $row2Values.Length

# Excel Range.Value expects a 2D array for a multi-cell range assignment.
$row2Array = New-Object 'object[,]' 1, $n
$row3Array = New-Object 'object[,]' 1, $n
for ($i = 0; $i -lt $n; $i++) {
    $row2Array[0, $i] = $row2Values[$i]
    $row3Array[0, $i] = $row3Values[$i]
}

# Write the two new rows into V2:AO2 and V3:AO3
$ws.Range("V2:AO2").Value = $row2Array
$ws.Range("V3:AO3").Value = $row3Array

# Match the updated view state left in the file: scrolled right a bit and
# a new active cell/selection further down in the newly-added columns.
$ws.Activate()
[void]$ws.Range("AD10").Select()
$excel.ActiveWindow.ScrollColumn = 12
